$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 454 (i.e. at 455/456), pushing the
# existing 455..482 block down to 457..484. Excel copies formatting
# (including the date-style on column D) from the row above automatically.
$ws.Rows("455:456").Insert()

# Row 455 - new weekly "Primera" quality record
$ws.Cells.Item(455, 1).Value2 = 9
$ws.Cells.Item(455, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(455, 3).Value2 = "Metropolitana"
$ws.Cells.Item(455, 4).Value2 = 45106
$ws.Cells.Item(455, 5).Value2 = 13
$ws.Cells.Item(455, 6).Value2 = 100112017
$ws.Cells.Item(455, 7).Value2 = "Apio"
$ws.Cells.Item(455, 8).Value2 = "Americana (o)"
$ws.Cells.Item(455, 9).Value2 = "Primera"
$ws.Cells.Item(455, 10).Value2 = 70
$ws.Cells.Item(455, 11).Value2 = 7000
$ws.Cells.Item(455, 12).Value2 = 8000
$ws.Cells.Item(455, 13).Value2 = 7500
$ws.Cells.Item(455, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(455, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(455, 16).Value2 = 1250
$ws.Cells.Item(455, 17).Value2 = 6
$ws.Cells.Item(455, 18).Value2 = "Hortaliza"

# Row 456 - new weekly "Segunda" quality record
$ws.Cells.Item(456, 1).Value2 = 9
$ws.Cells.Item(456, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(456, 3).Value2 = "Metropolitana"
$ws.Cells.Item(456, 4).Value2 = 45106
$ws.Cells.Item(456, 5).Value2 = 13
$ws.Cells.Item(456, 6).Value2 = 100112017
$ws.Cells.Item(456, 7).Value2 = "Apio"
$ws.Cells.Item(456, 8).Value2 = "Americana (o)"
$ws.Cells.Item(456, 9).Value2 = "Segunda"
$ws.Cells.Item(456, 10).Value2 = 52
$ws.Cells.Item(456, 11).Value2 = 6000
$ws.Cells.Item(456, 12).Value2 = 6000
$ws.Cells.Item(456, 13).Value2 = 6000
$ws.Cells.Item(456, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(456, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(456, 16).Value2 = 1000
$ws.Cells.Item(456, 17).Value2 = 6
$ws.Cells.Item(456, 18).Value2 = "Hortaliza"
